# Refactor: add an "interest rates" scenario sheet next to "metric overrides".

$wb = $excel.ActiveWorkbook
$wsOverrides = $wb.Worksheets.Item(1)

# --- Add the new "interest rates" sheet right after "metric overrides" ---
$ws = $wb.Worksheets.Add($null, $wsOverrides)
$ws.Name = "interest rates"

# --- Header / template rows -------------------------------------------------
$ws.Range("A1").Value = "Template"
$ws.Range("B1").Value = "InterestRates"

$ws.Range("A2").Value = "Date"
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "Type"
$ws.Range("D2").Value = "Tenor"
$ws.Range("E2").Value = "Maturity"
$ws.Range("F2").Value = "Rate"

# --- Data rows (Euribor curve as of 2022-12-31, serial 44926) ---------------
$rows = @(
    @{ Row=3;  Type="Spot"; Tenor="3m"; Maturity=$null; Rate=0.0305 },
    @{ Row=4;  Type="Spot"; Tenor="6m"; Maturity=$null; Rate=0.0295 },
    @{ Row=5;  Type="Zero"; Tenor=$null; Maturity="1m"; Rate=0.0318 },
    @{ Row=6;  Type="Zero"; Tenor=$null; Maturity="1y"; Rate=0.0286 },
    @{ Row=7;  Type="Zero"; Tenor=$null; Maturity="10y"; Rate=0.0255 },
    @{ Row=8;  Type="Zero"; Tenor=$null; Maturity="20y"; Rate=0.0265 },
    @{ Row=9;  Type="Zero"; Tenor=$null; Maturity="30y"; Rate=0.0270 },
    @{ Row=10; Type="Swap"; Tenor="3M"; Maturity="1Y"; Rate=0.0285 },
    @{ Row=11; Type="Swap"; Tenor="3M"; Maturity="10Y"; Rate=0.0255 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Range("A$row").Value = 44926
    $ws.Range("A$row").NumberFormat = "m/d/yy"
    $ws.Range("B$row").Value = "Euribor"
    $ws.Range("C$row").Value = $r.Type
    if ($r.Tenor) { $ws.Range("D$row").Value = $r.Tenor }
    if ($r.Maturity) { $ws.Range("E$row").Value = $r.Maturity }
    $ws.Range("F$row").Value = $r.Rate
    $ws.Range("F$row").NumberFormat = "0.00%"
}

# --- Update the "metric overrides" sheet's remembered selection -------------
$wsOverrides.Range("A2").Select()

# --- Cosmetics: column width + selection on the new sheet, then activate it -
$ws.Columns.Item(1).ColumnWidth = 14.90625
$ws.Range("A5").Select()
$ws.Activate()
